# update beauty, MG, PSS dashboard
# The first data row (row 2) on Sheet1 gets a freshly generated case id.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-A2735JP7"
